$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append AQI sensor rows 36-63 for "Loni, Ghaziabad - UPPCB" ---

# Pre-populate B/C values for the new rows
$ws.Range("A36").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B36").Value = "Thursday, 01 Apr 2021 10:00 AM"
$ws.Range("C36").Value = 246
$ws.Range("A37").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B37").Value = "Friday, 02 Apr 2021 10:00 AM"
$ws.Range("C37").Value = 290
$ws.Range("A38").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B38").Value = "Saturday, 03 Apr 2021 10:00 AM"
$ws.Range("C38").Value = 204
$ws.Range("A39").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B39").Value = "Sunday, 04 Apr 2021 10:00 AM"
$ws.Range("C39").Value = 200
$ws.Range("A40").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B40").Value = "Monday, 05 Apr 2021 10:00 AM"
$ws.Range("C40").Value = 257
$ws.Range("A41").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B41").Value = "Tuesday, 06 Apr 2021 10:00 AM"
$ws.Range("C41").Value = 342
$ws.Range("A42").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B42").Value = "Wednesday, 07 Apr 2021 10:00 AM"
$ws.Range("C42").Value = 336
$ws.Range("A43").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B43").Value = "Friday, 09 Apr 2021 10:00 AM"
$ws.Range("C43").Value = 206
$ws.Range("A44").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B44").Value = "Saturday, 10 Apr 2021 10:00 AM"
$ws.Range("C44").Value = 227
$ws.Range("A45").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B45").Value = "Sunday, 11 Apr 2021 10:00 AM"
$ws.Range("C45").Value = 275
$ws.Range("A46").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B46").Value = "Monday, 12 Apr 2021 10:00 AM"
$ws.Range("C46").Value = 325
$ws.Range("A47").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B47").Value = "Wednesday, 14 Apr 2021 10:00 AM"
$ws.Range("C47").Value = 312
$ws.Range("A48").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B48").Value = "Thursday, 15 Apr 2021 10:00 AM"
$ws.Range("C48").Value = 306
$ws.Range("A49").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B49").Value = "Friday, 16 Apr 2021 10:00 AM"
$ws.Range("C49").Value = 352
$ws.Range("A50").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B50").Value = "Saturday, 17 Apr 2021 10:00 AM"
$ws.Range("C50").Value = 169
$ws.Range("A51").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B51").Value = "Sunday, 18 Apr 2021 10:00 AM"
$ws.Range("C51").Value = 222
$ws.Range("A52").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B52").Value = "Monday, 19 Apr 2021 10:00 AM"
$ws.Range("C52").Value = 248
$ws.Range("A53").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B53").Value = "Tuesday, 20 Apr 2021 10:00 AM"
$ws.Range("C53").Value = 261
$ws.Range("A54").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B54").Value = "Wednesday, 21 Apr 2021 10:00 AM"
$ws.Range("C54").Value = 193
$ws.Range("A55").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B55").Value = "Thursday, 22 Apr 2021 10:00 AM"
$ws.Range("C55").Value = 168
$ws.Range("A56").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B56").Value = "Friday, 23 Apr 2021 10:00 AM"
$ws.Range("C56").Value = 183
$ws.Range("A57").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B57").Value = "Saturday, 24 Apr 2021 10:00 AM"
$ws.Range("C57").Value = 130
$ws.Range("A58").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B58").Value = "Sunday, 25 Apr 2021 10:00 AM"
$ws.Range("C58").Value = 227
$ws.Range("A59").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B59").Value = "Monday, 26 Apr 2021 10:00 AM"
$ws.Range("C59").Value = 308
$ws.Range("A60").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B60").Value = "Tuesday, 27 Apr 2021 10:00 AM"
$ws.Range("C60").Value = 364
$ws.Range("A61").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B61").Value = "Wednesday, 28 Apr 2021 10:00 AM"
$ws.Range("C61").Value = 412
$ws.Range("A62").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B62").Value = "Thursday, 29 Apr 2021 10:00 AM"
$ws.Range("C62").Value = 400
$ws.Range("A63").Value = "Loni, Ghaziabad - UPPCB"
$ws.Range("B63").Value = "Friday, 30 Apr 2021 10:00 AM"
$ws.Range("C63").Value = 369

# Copy the formatting (column styles) from the last existing data row (35)
# down across the newly added rows (36-63), matching the column-A/column-B styling
# used throughout the sheet (s="1" on A, s="2" on B).
$ws.Range("A35:C35").Copy()
$ws.Range("A36:C63").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the 15.75pt row height used by every other data row.
$ws.Range("A36:C63").RowHeight = 15.75

# Update the view: scroll so row 49 is at the top and select D57,
# mirroring the sheetView/selection recorded after the edit.
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D57").Select()

